$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grouping Info")

# Move the Price column (currently C) into column B, overwriting the
# Quantity values that used to live there, then drop the old Price (C)
# and Section (D) columns so only Item/Price remain.
$ws.Columns.Item(2).Delete()
$ws.Columns.Item(3).Delete()

# Update the selection to match the new active cell used after the edit.
$ws.Range("I19").Select()
